$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(12.64940783451405,8.87297407047606,17.29567131579251),
    @(17.19315151055925,14.92123702245383,19.84787293183696),
    @(16.40468380250744,14.19063040052704,18.79029149418443),
    @(20.59049357227114,18.12573770777175,23.46395483281349),
    @(18.56608644165994,14.88052166590214,22.53952347035616),
    @(11.90957674751932,8.52424745309159,15.22585451409978),
    @(4.471452986184406,3.036840082564717,6.449868774724746),
    @(17.40796613991894,14.12443497733679,20.66393268986968),
    @(4.567717066225484,3.095320806278899,6.465860688471674),
    @(10.0773205646123,6.879555347414891,14.22188615504995),
    @(16.7482139881463,12.39234942823371,21.58082079057044),
    @(5.982724749644271,4.027506516640751,8.82734072010642),
    @(7.79781817630193,5.332840851722391,11.27576955096275),
    @(23.85436960248326,20.55963034384466,27.52005792061857),
    @(18.53168151387659,14.85823739067681,22.50305600824254),
    @(16.65670924536424,12.34287376234699,21.93816666016205),
    @(16.40454908286583,14.19041250559461,18.78062833932014),
    @(10.15601724297252,6.884926687264382,14.36805766056095),
    @(20.31342931758741,16.18180761712782,24.70901019378979),
    @(22.22579780040033,18.72548215087375,26.31252602149752),
    @(4.553263822298934,3.047934610514206,6.649768317838993),
    @(25.39769808955764,22.08482759129611,28.90981938291566),
    @(18.54574366160662,14.90207883911431,22.54413709699854),
    @(16.27699509538217,13.43317746658813,19.17668112580869),
    @(22.1785222975604,18.21181979216581,25.9904949834773),
    @(5.989387198237481,4.028508269157857,8.568317858732545),
    @(12.52024249894082,8.777523670029485,16.40637402907041),
    @(26.62090122180022,22.96864924919585,30.66479253171784),
    @(16.31043172970823,11.92898266880038,21.42662595570505),
    @(7.740145351549794,4.958795120029158,11.22566776381867),
    @(5.934681732754369,4.197305479281908,8.356643188052139),
    @(7.799629219495676,5.345165361350537,11.2436967957312),
    @(13.35897716803587,9.196031523773636,18.29699223619431),
    @(19.46052718701372,15.36001043294635,24.02461989563266),
    @(11.92069088647269,8.439963927105936,15.29958730713738),
    @(23.62083559825562,19.66148300918483,27.62014361995718),
    @(7.74345236260048,5.245360899099869,11.142886998095),
    @(16.95146631264759,12.0042119145889,22.62356814905452),
    @(19.45788449061237,16.42130254448788,22.80708777984024),
    @(11.33709018521082,8.255474286292674,14.82440071664059),
    @(10.17515734920535,6.907178955328532,14.71995117391737),
    @(16.30390812015601,11.91673110100223,21.4098409141159),
    @(9.408510479868228,6.651938744917373,12.66739987031373),
    @(10.15681294786328,6.8839642811343,14.77859627329835),
    @(7.405950356889289,5.072076515531126,10.43364252143276),
    @(4.664981953953927,3.158881516763826,6.903308379129665),
    @(16.75062969232651,12.36071918295833,21.64130689429964),
    @(17.93898590554933,15.64973527648579,20.91905683695518),
    @(23.58173296197513,19.64000787398198,27.56957439655841),
    @(17.20155061189376,14.94101446170502,19.84810179302152),
    @(7.527480999789885,5.219016847866532,10.60073036612178),
    @(17.10149693854164,12.16549140195028,22.92415815216181),
    @(23.81801380284467,20.49509015518096,27.47115933393028),
    @(22.13670902297333,18.44514470164806,26.18413786734458),
    @(7.733572815701888,5.315034682462072,10.81714562733303),
    @(5.943388721477787,4.161115548954537,8.438714688947147),
    @(23.55908009918106,19.78990904115806,27.70514728553471),
    @(9.871728264364299,6.892176035494385,13.31583743143451),
    @(25.60979680933319,22.33262181890576,29.22139610134668),
    @(7.810652012461661,5.382070967492516,11.24357696266625),
    @(10.11906916637217,6.862944873195753,14.5711477763813),
    @(7.748072975526186,5.180294017335905,11.38775252071242),
    @(25.56199190842872,22.23061743245779,29.18062854621577),
    @(5.827008028310416,4.004718507272139,8.172336498434502),
    @(13.21261337999587,9.308045050200583,17.75406229718666),
    @(13.03461264572013,9.085399235147825,18.01442426948678),
    @(20.70239789602343,17.24350766228922,23.96907991885501),
    @(4.603576291968313,3.123438717600944,6.824201659949625),
    @(19.42701880176686,16.44403753869875,22.45681371855858),
    @(13.408320975414,10.50403286481913,16.19671092016774),
    @(14.9633103098299,11.28332997872225,18.70112265538329)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}
